$p = $ppt.ActivePresentation

# -----------------------------------------------------------------------
# Insert a brand-new slide at position 2 ("Deskripsi Umum Perangkat Lunak")
# using the same "Title and Content" layout (ppLayoutObject) that all the
# other content slides use. Everything that used to be slide 2..9 is
# pushed down to 3..10 automatically.
# -----------------------------------------------------------------------
$newSlide = $p.Slides.Add(2, 16)

# COM places the content placeholder first, the title placeholder second.
$content = $newSlide.Shapes.Item(1)
$title   = $newSlide.Shapes.Item(2)

# ---- Body / content placeholder -----------------------------------------
$bodyText = "Kehidupan manusia sekarang ini tidak dapat terlepas dari informasi, hal ini menjadi salah satu penyebab pentingnya keterlibatan bidang teknologi. Hasil pemanggilan ini akan di olah menjadi sebuah file keluaran (output) yang akan di tampilkan dalam bentuk gambar dan teks. Sehingga dengan digunakannya do’a-do’a sehari-hari yang sudah terkomputerisasi ini dapat meningkatkan pendapatan keimanan seseorang."

$tr = $content.TextFrame.TextRange
$tr.Text = $bodyText

# No bullet / no indent on this paragraph.
$tr.ParagraphFormat.Bullet.Visible = 0
$tr.ParagraphFormat.LeftIndent = 0
$tr.ParagraphFormat.FirstLineIndent = 0

# Italicize "file " and "output" (1-based character offsets).
$contentItalic1 = $tr.Characters(197, 5)
$contentItalic1.Font.Italic = -1
$contentItalic2 = $tr.Characters(212, 6)
$contentItalic2.Font.Italic = -1

# Trailing empty paragraph (matches the source deck which ends with a
# blank line after the description paragraph).
$content.TextFrame.TextRange.InsertAfter("`r")

# ---- Title placeholder ----------------------------------------------------
$titleTr = $title.TextFrame.TextRange
$titleTr.Text = "Deskripsi Umum Perangkat Lunak"
$titleTr.Font.Bold = -1

# Shrink-to-fit title text box (normAutofit), same as the rest of the deck.
$titleFrame2 = $title.TextFrame2
$titleFrame2.AutoSize = 2
$titleFrame2.FontScale = 90
